$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 45.71598933333333
$ws.Range("H2").Value = 137.147968
$ws.Range("I2").Value = 0.6549002937372808
$ws.Range("J2").Value = 0.6549002937372808
$ws.Range("O2").Value = 0.9347132976570145
$ws.Range("P2").Value = 0.9347132976570145
$ws.Range("Q2").Value = 406.0828051695431
$ws.Range("R2").Value = 3654.745246525888
$ws.Range("S2").Value = 0.6121440131957211
$ws.Range("T2").Value = 0.6121440131957211
# Row 3
$ws.Range("G3").Value = 45.71598933333333
$ws.Range("H3").Value = 137.147968
$ws.Range("I3").Value = 0.6549002937372808
$ws.Range("J3").Value = 0.6549002937372808
$ws.Range("M3").Value = 0.616144
$ws.Range("N3").Value = 1.848432
$ws.Range("O3").Value = 0.06483569448352988
$ws.Range("P3").Value = 0.0648356944835299
$ws.Range("Q3").Value = 28.16763253179733
$ws.Range("R3").Value = 253.508692786176
$ws.Range("S3").Value = 0.04246091536192431
$ws.Range("T3").Value = 0.04246091536192433
# Row 4
$ws.Range("G4").Value = 45.71598933333333
$ws.Range("H4").Value = 137.147968
$ws.Range("I4").Value = 0.6549002937372808
$ws.Range("J4").Value = 0.6549002937372808
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.004286
$ws.Range("N4").Value = 0.012858
$ws.Range("O4").Value = 0.0004510078594555965
$ws.Range("P4").Value = 0.0004510078594555965
$ws.Range("Q4").Value = 0.1959387302826667
$ws.Range("R4").Value = 1.763448572544
$ws.Range("S4").Value = 0.0002953651796352924
$ws.Range("T4").Value = 0.0002953651796352924
# Row 5
$ws.Range("I5").Value = 0.1818108415648851
$ws.Range("J5").Value = 0.1818108415648851
$ws.Range("O5").Value = 0.9347132976570145
$ws.Range("P5").Value = 0.9347132976570145
$ws.Range("S5").Value = 0.1699410112689107
$ws.Range("T5").Value = 0.1699410112689107
# Row 6
$ws.Range("I6").Value = 0.1818108415648851
$ws.Range("J6").Value = 0.1818108415648851
$ws.Range("M6").Value = 0.616144
$ws.Range("N6").Value = 1.848432
$ws.Range("O6").Value = 0.06483569448352988
$ws.Range("P6").Value = 0.0648356944835299
$ws.Range("Q6").Value = 7.819787262992
$ws.Range("R6").Value = 70.37808536692799
$ws.Range("S6").Value = 0.01178783217749435
$ws.Range("T6").Value = 0.01178783217749435
# Row 7
$ws.Range("I7").Value = 0.1818108415648851
$ws.Range("J7").Value = 0.1818108415648851
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.004286
$ws.Range("N7").Value = 0.012858
$ws.Range("O7").Value = 0.0004510078594555965
$ws.Range("P7").Value = 0.0004510078594555965
$ws.Range("Q7").Value = 0.05439573899799999
$ws.Range("R7").Value = 0.4895616509819999
$ws.Range("S7").Value = 8.199811847999941 / 100000
$ws.Range("T7").Value = 8.199811847999941 / 100000
# Row 8
$ws.Range("G8").Value = 11.24784666666667
$ws.Range("H8").Value = 33.74354
$ws.Range("I8").Value = 0.161130015850732
$ws.Range("J8").Value = 0.161130015850732
$ws.Range("O8").Value = 0.9347132976570145
$ws.Range("P8").Value = 0.9347132976570145
$ws.Range("Q8").Value = 99.91158877068221
$ws.Range("R8").Value = 899.2042989361398
$ws.Range("S8").Value = 0.1506103684673647
$ws.Range("T8").Value = 0.1506103684673647
# Row 9
$ws.Range("G9").Value = 11.24784666666667
$ws.Range("H9").Value = 33.74354
$ws.Range("I9").Value = 0.161130015850732
$ws.Range("J9").Value = 0.161130015850732
$ws.Range("M9").Value = 0.616144
$ws.Range("N9").Value = 1.848432
$ws.Range("O9").Value = 0.06483569448352988
$ws.Range("P9").Value = 0.0648356944835299
$ws.Range("Q9").Value = 6.930293236586667
$ws.Range("R9").Value = 62.37263912928
$ws.Range("S9").Value = 0.01044697647982439
$ws.Range("T9").Value = 0.01044697647982439
# Row 10
$ws.Range("G10").Value = 11.24784666666667
$ws.Range("H10").Value = 33.74354
$ws.Range("I10").Value = 0.161130015850732
$ws.Range("J10").Value = 0.161130015850732
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.004286
$ws.Range("N10").Value = 0.012858
$ws.Range("O10").Value = 0.0004510078594555965
$ws.Range("P10").Value = 0.0004510078594555965
$ws.Range("Q10").Value = 0.04820827081333333
$ws.Range("R10").Value = 0.4338744373199999
$ws.Range("S10").Value = 7.267090354288495 / 100000
$ws.Range("T10").Value = 7.267090354288497 / 100000
# Row 11
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.1507006666666667
$ws.Range("H11").Value = 0.452102
$ws.Range("I11").Value = 0.00215884884710222
$ws.Range("J11").Value = 0.00215884884710222
$ws.Range("O11").Value = 0.9347132976570145
$ws.Range("P11").Value = 0.9347132976570145
$ws.Range("Q11").Value = 1.338633383053556
$ws.Range("R11").Value = 12.047700447482
$ws.Range("S11").Value = 0.00201790472501796
$ws.Range("T11").Value = 0.00201790472501796
# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.1507006666666667
$ws.Range("H12").Value = 0.452102
$ws.Range("I12").Value = 0.00215884884710222
$ws.Range("J12").Value = 0.00215884884710222
$ws.Range("M12").Value = 0.616144
$ws.Range("N12").Value = 1.848432
$ws.Range("O12").Value = 0.06483569448352988
$ws.Range("P12").Value = 0.0648356944835299
$ws.Range("Q12").Value = 0.09285331156266667
$ws.Range("R12").Value = 0.835679804064
$ws.Range("S12").Value = 0.0001399704642868403
$ws.Range("T12").Value = 0.0001399704642868403
# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.1507006666666667
$ws.Range("H13").Value = 0.452102
$ws.Range("I13").Value = 0.00215884884710222
$ws.Range("J13").Value = 0.00215884884710222
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.004286
$ws.Range("N13").Value = 0.012858
$ws.Range("O13").Value = 0.0004510078594555965
$ws.Range("P13").Value = 0.0004510078594555965
$ws.Range("Q13").Value = 0.0006459030573333334
$ws.Range("R13").Value = 0.005813127516
$ws.Range("S13").Value = 9.736577974197546 / 10000000
$ws.Range("T13").Value = 9.736577974197546 / 10000000
